$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.979.11'
$ws.Range("E2").Value = '  -1.69%  '

$ws.Range("D3").Value = '1.910.90'
$ws.Range("E3").Value = '  -2.82%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.94'
$ws.Range("E5").Value = '  -0.46%  '

$ws.Range("E6").Value = '  -0.24%  '

$ws.Range("E7").Value = '  -1.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3823'
$ws.Range("E8").Value = '  -1.76%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07719'
$ws.Range("E9").Value = '  -2.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9807'
$ws.Range("E10").Value = '  -0.71%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.12'
$ws.Range("E11").Value = '  -2.72%  '

$ws.Range("D12").Value = '1.884.81'
$ws.Range("E12").Value = '  -4.66%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.671'
$ws.Range("E13").Value = '  -1.93%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.940'
$ws.Range("E14").Value = '  -2.72%  '

$ws.Range("E15").Value = '  -0.69%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").Value = '  -0.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '83.97'
$ws.Range("E17").Value = '  -3.86%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009464'
$ws.Range("E18").Value = '  -4.63%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.74'
$ws.Range("E19").Value = '  -2.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.34%  '

$ws.Range("D21").Value = '28.949.39'
$ws.Range("E21").Value = '  -1.86%  '

$ws.Range("E22").Value = '  -3.64%  '

$ws.Range("E23").Value = '  -1.89%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.086'
$ws.Range("E24").Value = '  -0.88%  '

$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.54'
$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.02'
$ws.Range("E26").Value = '  -2.21%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.673'
$ws.Range("E27").Value = '  -1.57%  '

$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '117.57'
$ws.Range("E28").Value = '  -1.45%  '

$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.853'
$ws.Range("E29").Value = '  -1.42%  '

$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09293'
$ws.Range("E30").Value = '  -1.20%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.8669'
$ws.Range("E31").Value = '  -0.17%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.074'
$ws.Range("E32").Value = '  -2.51%  '

$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.252'
$ws.Range("E33").Value = '  -4.75%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.022'
$ws.Range("E34").Value = '  -3.55%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05737'
$ws.Range("E35").Value = '  -0.80%  '

$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.154'
$ws.Range("E36").Value = '  -0.30%  '

$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.001'
$ws.Range("E37").Value = '  -0.24%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02044'
$ws.Range("E38").Value = '  -2.28%  '

$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5515'
$ws.Range("E39").Value = '  -3.13%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.415'
$ws.Range("E40").Value = '  -3.41%  '

$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1758'
$ws.Range("E41").Value = '  -1.87%  '

$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.841'
$ws.Range("E42").Value = '  +3.27%  '

$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.323'
$ws.Range("E43").Value = '  -2.79%  '

$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5188'
$ws.Range("E44").Value = '  -2.55%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.25'
$ws.Range("E45").Value = '  -3.82%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06858'
$ws.Range("E46").Value = '  -0.76%  '

$ws.Range("B47").Value = 'PEPE'
$ws.Range("C47").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000002623'
$ws.Range("E47").Value = '  -5.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.061'
$ws.Range("E48").Value = '  -3.49%  '

$ws.Range("B49").Value = 'Quant'
$ws.Range("C49").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '111.00'
$ws.Range("E49").Value = '  -1.09%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.782'
$ws.Range("E50").Value = '  -2.13%  '

$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("E51").Value = '  -0.14%  '
